$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 93-95
$ws.Range("B93").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"
$ws.Range("B94").Value = "['BTCUSD.SPOT']"
$ws.Range("B95").Value = "['BTCUSD.SPOT']"

# Add new rows 96-100 (dates written as text, mirroring the existing text-date column)
$ws.Range("A96").Value = "'2025-09-15"
$ws.Range("A96").ClearFormats()
$ws.Range("B96").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"

$ws.Range("A97").Value = "'2025-09-16"
$ws.Range("A97").ClearFormats()
$ws.Range("B97").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"

$ws.Range("A98").Value = "'2025-09-17"
$ws.Range("A98").ClearFormats()
$ws.Range("B98").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"

$ws.Range("A99").Value = "'2025-09-18"
$ws.Range("A99").ClearFormats()
$ws.Range("B99").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']"

$ws.Range("A100").Value = "'2025-09-19"
$ws.Range("A100").ClearFormats()
$ws.Range("B100").Value = "['USD.SOFR.CSA_USD']"
